# Generate Report for Archive
# - Update localization status text from "Ready for handoff" to "In Translation"
#   on the Overview sheet (zh-cn / de-de status columns) and on each
#   per-locale detail sheet (Status column).
# - Refresh the "Status" column widths so they reflect the new (shorter)
#   text, as AutoFit would after the report content changed.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "In Translation"

# Overview sheet: column E = zh-cn status, column F = de-de status
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus

# Per-locale detail sheets: column C = Status
$wsZhCn.Range("C2").Value = $newStatus
$wsDeDe.Range("C2").Value = $newStatus

# Resize the status columns to match the new text width.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
